# Loan Provisioning 10 test cases
#
# 1) Rotate the names of the three "Acc_Upfront*" sheets:
#      Acc_Upfront1 -> Acc_Upfront2
#      Acc_Upfront2 -> Acc_Upfront3
#      Acc_Upfront3 -> Acc_Upfront1
#    (use a temporary name to avoid collisions while rotating)
#
# 2) Update the remembered selection on a couple of sheets, and move the
#    active-sheet / active-tab onto the sheet that is now named
#    "Acc_Upfront2" (the former "Acc_Upfront1").

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Acc_Upfront1").Name = "Acc_UpfrontTemp"
$wb.Worksheets.Item("Acc_Upfront3").Name = "Acc_Upfront1"
$wb.Worksheets.Item("Acc_Upfront2").Name = "Acc_Upfront3"
$wb.Worksheets.Item("Acc_UpfrontTemp").Name = "Acc_Upfront2"

# Transactions sheet: move the remembered selection from C6 to D2.
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("D2").Select()

# Acc_Upfront1 (was Acc_Upfront3): keep selection at F3, drop it as the
# active tab.
$wsUpfront1 = $wb.Worksheets.Item("Acc_Upfront1")
$wsUpfront1.Activate()
$wsUpfront1.Range("F3").Select()

# Acc_Upfront2 (was Acc_Upfront1): becomes the active/selected tab, with
# the remembered selection moved from E2 to F22.
$wsUpfront2 = $wb.Worksheets.Item("Acc_Upfront2")
$wsUpfront2.Activate()
$wsUpfront2.Range("F22").Select()
